$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = 2
    4  = -1
    5  = -2
    6  = -3
    7  = 2
    8  = 2
    10 = -2
    12 = -3
    13 = -4
    14 = -6
    15 = -2
    16 = 11
    17 = 2
    18 = 1
    20 = 3
    21 = -3
    22 = -3
    23 = 1
    24 = -3
    25 = 1
    26 = -4
    28 = -1
    29 = -6
    30 = 2
    31 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
